# Update cached Universalis market-price / leve-profit figures across the
# per-job "Marilith Profits" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# These are raw scraped price snapshots (no formulas), so the refreshed
# numbers are applied as straight value writes, one cell at a time.

$wb = $excel.ActiveWorkbook

# Every cell whose cached value changed (or was newly populated) in this
# refresh, keyed by worksheet name + A1 reference.
$updates = @(
    @{ Sheet = "ALC"; Cell = "H64"; Value = 4999.9 },
    @{ Sheet = "ALC"; Cell = "I64"; Value = 4999.9 },
    @{ Sheet = "ALC"; Cell = "K64"; Value = 4999.9 },
    @{ Sheet = "ALC"; Cell = "M64"; Value = -4751.9 },
    @{ Sheet = "ALC"; Cell = "H67"; Value = 4999.9 },
    @{ Sheet = "ALC"; Cell = "I67"; Value = 4999.9 },
    @{ Sheet = "ALC"; Cell = "K67"; Value = 4999.9 },
    @{ Sheet = "ALC"; Cell = "M67"; Value = -4141.9 },
    @{ Sheet = "ALC"; Cell = "H94"; Value = 18066.666 },
    @{ Sheet = "ALC"; Cell = "I94"; Value = 18066.666 },
    @{ Sheet = "ALC"; Cell = "K94"; Value = 18066.666 },
    @{ Sheet = "ALC"; Cell = "M94"; Value = -17615.666 },
    @{ Sheet = "ALC"; Cell = "H98"; Value = 849.6 },
    @{ Sheet = "ALC"; Cell = "I98"; Value = 493.66666 },
    @{ Sheet = "ALC"; Cell = "J98"; Value = 4053 },
    @{ Sheet = "ALC"; Cell = "K98"; Value = 493.66666 },
    @{ Sheet = "ALC"; Cell = "L98"; Value = 4053 },
    @{ Sheet = "ALC"; Cell = "M98"; Value = 1004.33334 },
    @{ Sheet = "ALC"; Cell = "N98"; Value = -7049 },
    @{ Sheet = "ALC"; Cell = "H107"; Value = 2663.9092 },
    @{ Sheet = "ALC"; Cell = "I107"; Value = 686.5 },
    @{ Sheet = "ALC"; Cell = "K107"; Value = 686.5 },
    @{ Sheet = "ALC"; Cell = "M107"; Value = 1233.5 },
    @{ Sheet = "ALC"; Cell = "H116"; Value = 5997.5 },
    @{ Sheet = "ALC"; Cell = "I116"; Value = 5995 },
    @{ Sheet = "ALC"; Cell = "K116"; Value = 5995 },
    @{ Sheet = "ALC"; Cell = "M116"; Value = -2553 },
    @{ Sheet = "ALC"; Cell = "H122"; Value = 849.6 },
    @{ Sheet = "ALC"; Cell = "I122"; Value = 493.66666 },
    @{ Sheet = "ALC"; Cell = "J122"; Value = 4053 },
    @{ Sheet = "ALC"; Cell = "K122"; Value = 1480.99998 },
    @{ Sheet = "ALC"; Cell = "L122"; Value = 12159 },
    @{ Sheet = "ALC"; Cell = "M122"; Value = 969.0000199999999 },
    @{ Sheet = "ALC"; Cell = "N122"; Value = -17059 },
    @{ Sheet = "ALC"; Cell = "H132"; Value = 2411.5625 },
    @{ Sheet = "ALC"; Cell = "J132"; Value = 0 },
    @{ Sheet = "ALC"; Cell = "L132"; Value = 0 },
    @{ Sheet = "ARM"; Cell = "H2"; Value = 1892.5 },
    @{ Sheet = "ARM"; Cell = "I2"; Value = 1020 },
    @{ Sheet = "ARM"; Cell = "K2"; Value = 1020 },
    @{ Sheet = "ARM"; Cell = "M2"; Value = -907 },
    @{ Sheet = "ARM"; Cell = "H116"; Value = 1892.5 },
    @{ Sheet = "ARM"; Cell = "I116"; Value = 1020 },
    @{ Sheet = "ARM"; Cell = "K116"; Value = 1020 },
    @{ Sheet = "ARM"; Cell = "M116"; Value = 1274 },
    @{ Sheet = "ARM"; Cell = "H132"; Value = 2078.8 },
    @{ Sheet = "ARM"; Cell = "I132"; Value = 997.5 },
    @{ Sheet = "ARM"; Cell = "K132"; Value = 2992.5 },
    @{ Sheet = "ARM"; Cell = "M132"; Value = -462.5 },
    @{ Sheet = "BSM"; Cell = "H3"; Value = 1892.5 },
    @{ Sheet = "BSM"; Cell = "I3"; Value = 1020 },
    @{ Sheet = "BSM"; Cell = "K3"; Value = 1020 },
    @{ Sheet = "BSM"; Cell = "M3"; Value = -906 },
    @{ Sheet = "BSM"; Cell = "H64"; Value = 954.94116 },
    @{ Sheet = "BSM"; Cell = "I64"; Value = 1138 },
    @{ Sheet = "BSM"; Cell = "J64"; Value = 855.0909 },
    @{ Sheet = "BSM"; Cell = "K64"; Value = 1138 },
    @{ Sheet = "BSM"; Cell = "L64"; Value = 855.0909 },
    @{ Sheet = "BSM"; Cell = "M64"; Value = -913 },
    @{ Sheet = "BSM"; Cell = "N64"; Value = -1305.0909 },
    @{ Sheet = "BSM"; Cell = "H67"; Value = 954.94116 },
    @{ Sheet = "BSM"; Cell = "I67"; Value = 1138 },
    @{ Sheet = "BSM"; Cell = "J67"; Value = 855.0909 },
    @{ Sheet = "BSM"; Cell = "K67"; Value = 1138 },
    @{ Sheet = "BSM"; Cell = "L67"; Value = 855.0909 },
    @{ Sheet = "BSM"; Cell = "M67"; Value = -358 },
    @{ Sheet = "BSM"; Cell = "N67"; Value = -2415.0909 },
    @{ Sheet = "BSM"; Cell = "H94"; Value = 2343.2222 },
    @{ Sheet = "BSM"; Cell = "I94"; Value = 2386.125 },
    @{ Sheet = "BSM"; Cell = "K94"; Value = 2386.125 },
    @{ Sheet = "BSM"; Cell = "M94"; Value = -1935.125 },
    @{ Sheet = "BSM"; Cell = "H99"; Value = 1463.3334 },
    @{ Sheet = "BSM"; Cell = "I99"; Value = 1533.1578 },
    @{ Sheet = "BSM"; Cell = "K99"; Value = 1533.1578 },
    @{ Sheet = "BSM"; Cell = "M99"; Value = -35.15779999999995 },
    @{ Sheet = "BSM"; Cell = "H107"; Value = 867.8182 },
    @{ Sheet = "BSM"; Cell = "I107"; Value = 805.6667 },
    @{ Sheet = "BSM"; Cell = "K107"; Value = 805.6667 },
    @{ Sheet = "BSM"; Cell = "M107"; Value = 1114.3333 },
    @{ Sheet = "CRP"; Cell = "H31"; Value = 3593.4211 },
    @{ Sheet = "CRP"; Cell = "I31"; Value = 2435.1667 },
    @{ Sheet = "CRP"; Cell = "K31"; Value = 2435.1667 },
    @{ Sheet = "CRP"; Cell = "M31"; Value = -2140.1667 },
    @{ Sheet = "CRP"; Cell = "H34"; Value = 3593.4211 },
    @{ Sheet = "CRP"; Cell = "I34"; Value = 2435.1667 },
    @{ Sheet = "CRP"; Cell = "K34"; Value = 2435.1667 },
    @{ Sheet = "CRP"; Cell = "M34"; Value = -2233.1667 },
    @{ Sheet = "CRP"; Cell = "H94"; Value = 1663.4286 },
    @{ Sheet = "CRP"; Cell = "I94"; Value = 1657.3334 },
    @{ Sheet = "CRP"; Cell = "K94"; Value = 1657.3334 },
    @{ Sheet = "CRP"; Cell = "M94"; Value = -1206.3334 },
    @{ Sheet = "CRP"; Cell = "H122"; Value = 756.2857 },
    @{ Sheet = "CRP"; Cell = "I122"; Value = 740.4 },
    @{ Sheet = "CRP"; Cell = "J122"; Value = 796 },
    @{ Sheet = "CRP"; Cell = "K122"; Value = 2221.2 },
    @{ Sheet = "CRP"; Cell = "L122"; Value = 2388 },
    @{ Sheet = "CRP"; Cell = "M122"; Value = 228.8000000000002 },
    @{ Sheet = "CRP"; Cell = "N122"; Value = -7288 },
    @{ Sheet = "CRP"; Cell = "H134"; Value = 2607.84 },
    @{ Sheet = "CRP"; Cell = "I134"; Value = 2488.4 },
    @{ Sheet = "CRP"; Cell = "J134"; Value = 3085.6 },
    @{ Sheet = "CRP"; Cell = "K134"; Value = 7465.200000000001 },
    @{ Sheet = "CRP"; Cell = "L134"; Value = 9256.799999999999 },
    @{ Sheet = "CRP"; Cell = "M134"; Value = -4930.200000000001 },
    @{ Sheet = "CRP"; Cell = "N134"; Value = -14326.8 },
    @{ Sheet = "CUL"; Cell = "H117"; Value = 833.6667 },
    @{ Sheet = "CUL"; Cell = "I117"; Value = 763 },
    @{ Sheet = "CUL"; Cell = "K117"; Value = 2289 },
    @{ Sheet = "CUL"; Cell = "M117"; Value = 1153 },
    @{ Sheet = "CUL"; Cell = "H122"; Value = 710.35297 },
    @{ Sheet = "CUL"; Cell = "I122"; Value = 879.8 },
    @{ Sheet = "CUL"; Cell = "J122"; Value = 639.75 },
    @{ Sheet = "CUL"; Cell = "K122"; Value = 7918.2 },
    @{ Sheet = "CUL"; Cell = "L122"; Value = 5757.75 },
    @{ Sheet = "CUL"; Cell = "M122"; Value = -5468.2 },
    @{ Sheet = "CUL"; Cell = "N122"; Value = -10657.75 },
    @{ Sheet = "GSM"; Cell = "H5"; Value = 41239.5 },
    @{ Sheet = "GSM"; Cell = "I5"; Value = 41239.5 },
    @{ Sheet = "GSM"; Cell = "K5"; Value = 41239.5 },
    @{ Sheet = "GSM"; Cell = "M5"; Value = -41127.5 },
    @{ Sheet = "GSM"; Cell = "H97"; Value = 536.5 },
    @{ Sheet = "GSM"; Cell = "I97"; Value = 367.58334 },
    @{ Sheet = "GSM"; Cell = "K97"; Value = 367.58334 },
    @{ Sheet = "GSM"; Cell = "M97"; Value = 128.41666 },
    @{ Sheet = "GSM"; Cell = "H113"; Value = 4550.5557 },
    @{ Sheet = "GSM"; Cell = "J113"; Value = 4798.5 },
    @{ Sheet = "GSM"; Cell = "L113"; Value = 4798.5 },
    @{ Sheet = "GSM"; Cell = "N113"; Value = -9138.5 },
    @{ Sheet = "LTW"; Cell = "H16"; Value = 359 },
    @{ Sheet = "LTW"; Cell = "I16"; Value = 359 },
    @{ Sheet = "LTW"; Cell = "K16"; Value = 359 },
    @{ Sheet = "LTW"; Cell = "M16"; Value = -189 },
    @{ Sheet = "LTW"; Cell = "H40"; Value = 2665.3333 },
    @{ Sheet = "LTW"; Cell = "I40"; Value = 2665.3333 },
    @{ Sheet = "LTW"; Cell = "K40"; Value = 2665.3333 },
    @{ Sheet = "LTW"; Cell = "M40"; Value = -2529.3333 },
    @{ Sheet = "LTW"; Cell = "H58"; Value = 1914 },
    @{ Sheet = "LTW"; Cell = "I58"; Value = 1914 },
    @{ Sheet = "LTW"; Cell = "K58"; Value = 1914 },
    @{ Sheet = "LTW"; Cell = "M58"; Value = -1654 },
    @{ Sheet = "LTW"; Cell = "H61"; Value = 886.5 },
    @{ Sheet = "LTW"; Cell = "I61"; Value = 249.5 },
    @{ Sheet = "LTW"; Cell = "K61"; Value = 249.5 },
    @{ Sheet = "LTW"; Cell = "M61"; Value = -47.5 },
    @{ Sheet = "LTW"; Cell = "H93"; Value = 400 },
    @{ Sheet = "LTW"; Cell = "I93"; Value = 400 },
    @{ Sheet = "LTW"; Cell = "K93"; Value = 400 },
    @{ Sheet = "LTW"; Cell = "M93"; Value = 848 },
    @{ Sheet = "LTW"; Cell = "H113"; Value = 886.5 },
    @{ Sheet = "LTW"; Cell = "I113"; Value = 249.5 },
    @{ Sheet = "LTW"; Cell = "K113"; Value = 249.5 },
    @{ Sheet = "LTW"; Cell = "M113"; Value = 1920.5 },
    @{ Sheet = "LTW"; Cell = "H122"; Value = 3670.8333 },
    @{ Sheet = "LTW"; Cell = "I122"; Value = 3504 },
    @{ Sheet = "LTW"; Cell = "K122"; Value = 10512 },
    @{ Sheet = "LTW"; Cell = "M122"; Value = -8062 },
    @{ Sheet = "LTW"; Cell = "H132"; Value = 5969.3076 },
    @{ Sheet = "LTW"; Cell = "I132"; Value = 4450.125 },
    @{ Sheet = "LTW"; Cell = "K132"; Value = 13350.375 },
    @{ Sheet = "LTW"; Cell = "M132"; Value = -10820.375 },
    @{ Sheet = "WVR"; Cell = "H62"; Value = 4542.857 },
    @{ Sheet = "WVR"; Cell = "I62"; Value = 4466.1665 },
    @{ Sheet = "WVR"; Cell = "K62"; Value = 4466.1665 },
    @{ Sheet = "WVR"; Cell = "M62"; Value = -3842.1665 },
    @{ Sheet = "WVR"; Cell = "H65"; Value = 4542.857 },
    @{ Sheet = "WVR"; Cell = "I65"; Value = 4466.1665 },
    @{ Sheet = "WVR"; Cell = "K65"; Value = 22330.8325 },
    @{ Sheet = "WVR"; Cell = "M65"; Value = -19210.8325 },
    @{ Sheet = "WVR"; Cell = "H107"; Value = 397.2 },
    @{ Sheet = "WVR"; Cell = "I107"; Value = 328.66666 },
    @{ Sheet = "WVR"; Cell = "K107"; Value = 985.9999799999999 },
    @{ Sheet = "WVR"; Cell = "M107"; Value = 934.0000200000001 },
    @{ Sheet = "WVR"; Cell = "H113"; Value = 5753.55 },
    @{ Sheet = "WVR"; Cell = "I113"; Value = 7792.2856 },
    @{ Sheet = "WVR"; Cell = "K113"; Value = 23376.8568 },
    @{ Sheet = "WVR"; Cell = "M113"; Value = -21206.8568 },
    @{ Sheet = "WVR"; Cell = "H122"; Value = 1807.375 },
    @{ Sheet = "WVR"; Cell = "I122"; Value = 1555.8 },
    @{ Sheet = "WVR"; Cell = "J122"; Value = 2226.6667 },
    @{ Sheet = "WVR"; Cell = "K122"; Value = 4667.4 },
    @{ Sheet = "WVR"; Cell = "L122"; Value = 6680.000100000001 },
    @{ Sheet = "WVR"; Cell = "M122"; Value = -2217.4 },
    @{ Sheet = "WVR"; Cell = "N122"; Value = -11580.0001 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}

# Cells whose cached value is no longer produced by this refresh (the
# upstream data source stopped reporting a figure for that column) get
# cleared outright rather than left stale.
$clears = @(
    @{ Sheet = "ALC"; Cell = "N132" }
)

foreach ($c in $clears) {
    $ws = $wb.Worksheets.Item($c.Sheet)
    $ws.Range($c.Cell).ClearContents()
}
